# Resolve merge conflict on the employees workbook:
#  - Fill in a new employee row (Maria Dutu) on the Employees sheet.
#  - Make "Employees" the active/selected sheet again (was "PracticeForm").
#  - Update the remembered selections on both of those sheets.

$wb = $excel.ActiveWorkbook

$wsEmployees    = $wb.Worksheets.Item("Employees")
$wsPracticeForm = $wb.Worksheets.Item("PracticeForm")

# New employee data on row 13 of the Employees sheet.
$wsEmployees.Range("A13").Value = "Maria"
$wsEmployees.Range("B13").Value = "Dutu"
$wsEmployees.Range("C13").Value = "Bucuresti"
$wsEmployees.Range("D13").Value = "mariadutu.d@gmail.com"
$wsEmployees.Range("F13").Value = "mihaela.stanciu@gmail.com"
$wsEmployees.Range("G13").Value = "C,Prolog,Haskell,Python,"
$wsEmployees.Range("H13").Value = "ux design,ui design"
$wsEmployees.Range("I13").Value = "ui designer"

# Restore PracticeForm's remembered selection first (while it is still the
# active sheet), then switch activation/selection back to Employees so it
# ends up as the tab that is selected when the workbook re-opens.
$wsPracticeForm.Range("C12").Select()

$wsEmployees.Activate()
$wsEmployees.Range("C18").Select()
